$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header/data columns onto the two new
# columns (G, H) before writing their values, so G1/H1 match the other
# header cells (style s=1) and G2/H2 match the other data cells (style s=3).
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("E2").Copy() | Out-Null
$ws.Range("G2:H2").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# New header row values
$ws.Range("G1").Value = "expectedThanksHeader"
$ws.Range("H1").Value = "expectedThanksBody"

# New data row values
$ws.Range("G2").Value = "Thank you for your order!"
$ws.Range("H2").Value = "Your order has been dispatched, and will arrive just as fast as the pony can get there!"

$ws.Range("H1").Select() | Out-Null
